# Skills Assessment_Jnr Dev.docx - add review comments on the
# "Include unit tests...", "Basic error handling" and
# "Ability to write simple tests (bonus)" bullet items.
#
# The existing document already has a comment (w:id="8") anchored on
# "Basic error handling". The target revision keeps that comment's text
# but re-anchors the numbering so that:
#   - a brand new comment ("Unit tests added") lands on the earlier
#     "Include unit tests for one or more endpoints" bullet and takes
#     comment id 8,
#   - the pre-existing "Status codes..." comment text is preserved on
#     "Basic error handling" but ends up as comment id 9,
#   - a brand new comment ("Complete?") lands on the later
#     "Ability to write simple tests (bonus)" bullet and takes id 10.
#
# Word/the COM object model always hands out the next free integer id
# when Comments.Add is called (lowest unused id), regardless of where
# in the document the anchor range sits. So to land the new "Include
# unit tests" comment on id 8 we must first free up id 8 by removing
# the old comment, then re-create it (with identical text/author) on
# the same anchor so it naturally becomes id 9, and only then add the
# third, brand-new comment (which becomes id 10).

$d = $word.ActiveDocument

function Get-RangeForText($text) {
    $r = $d.Content.Duplicate
    $r.Find.ClearFormatting()
    $r.Find.Text = $text
    $r.Find.Execute() | Out-Null
    return $r
}

function Get-CommentByScopeText($text) {
    for ($i = 1; $i -le $d.Comments.Count; $i++) {
        $candidate = $d.Comments($i)
        if ($candidate.Scope.Text -eq $text) {
            return $candidate
        }
    }
    return $null
}

# Locate + remember the range of the existing "Basic error handling"
# comment before touching anything so we can re-anchor a replacement
# comment on exactly the same text after the old one is removed.
$basicErrorHandlingComment = Get-CommentByScopeText("Basic error handling")
$basicErrorHandlingRange = Get-RangeForText("Basic error handling")
$preservedAuthor = $basicErrorHandlingComment.Author
$preservedInitial = $basicErrorHandlingComment.Initial
$preservedText = $basicErrorHandlingComment.Range.Text

# Free up its id (8) by deleting it; we'll recreate identical content
# on the same anchor once the new "Include unit tests" comment has
# claimed id 8.
$basicErrorHandlingComment.Delete()

# New comment #1: "Include unit tests for one or more endpoints" ->
# becomes comment id 8 (lowest free id at this point).
$unitTestsRange = Get-RangeForText("Include unit tests for one or more endpoints")
$unitTestsComment = $d.Comments.Add($unitTestsRange, "Unit tests added")
$unitTestsComment.Author = $preservedAuthor
$unitTestsComment.Initial = $preservedInitial

# Re-create the original "Basic error handling" comment with its
# original text/author -> becomes comment id 9.
$basicErrorHandlingComment2 = $d.Comments.Add($basicErrorHandlingRange, $preservedText)
$basicErrorHandlingComment2.Author = $preservedAuthor
$basicErrorHandlingComment2.Initial = $preservedInitial

# New comment #2: "Ability to write simple tests (bonus)" -> becomes
# comment id 10.
$bonusTestsRange = Get-RangeForText("Ability to write simple tests (bonus)")
$bonusTestsComment = $d.Comments.Add($bonusTestsRange, "Complete?")
$bonusTestsComment.Author = $preservedAuthor
$bonusTestsComment.Initial = $preservedInitial

Write-Host "Comments now:" $d.Comments.Count
